# Refactor test endpoint for customer put over billing address
# Update the Billing sheet's address_1 value (E2) and move the active
# selection to E3, as captured in the latest workbook save state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Billing")

# Update the billing address_1 cell value
$ws.Range("E2").Value = "1 Aeropost Way"

# Activate the sheet and move selection to E3 to match saved view state
$ws.Activate()
$ws.Range("E3").Select()
